# Update the 100 arithmetic answer cells in the table to match the
# regenerated "answers-of-addition_and_subtraction_within_100" output.
# Each old equation string is unique in the document, so a simple
# whole-document Find/Replace (whole word match, no wildcards) for each
# pair is sufficient and unambiguous.
$d = $word.ActiveDocument

$d.Content.Find.Execute("40+17=57", $true, $false, $false, $false, $false, $true, 1, $false, "45+12=57", 2) | Out-Null
$d.Content.Find.Execute("55-48=7", $true, $false, $false, $false, $false, $true, 1, $false, "87-47=40", 2) | Out-Null
$d.Content.Find.Execute("93-29=64", $true, $false, $false, $false, $false, $true, 1, $false, "54+21=75", 2) | Out-Null
$d.Content.Find.Execute("9+83=92", $true, $false, $false, $false, $false, $true, 1, $false, "45-43=2", 2) | Out-Null
$d.Content.Find.Execute("35+48=83", $true, $false, $false, $false, $false, $true, 1, $false, "84-84=0", 2) | Out-Null
$d.Content.Find.Execute("31+24=55", $true, $false, $false, $false, $false, $true, 1, $false, "92-89=3", 2) | Out-Null
$d.Content.Find.Execute("1+54=55", $true, $false, $false, $false, $false, $true, 1, $false, "73-36=37", 2) | Out-Null
$d.Content.Find.Execute("45-34=11", $true, $false, $false, $false, $false, $true, 1, $false, "49+5=54", 2) | Out-Null
$d.Content.Find.Execute("46+43=89", $true, $false, $false, $false, $false, $true, 1, $false, "60+32=92", 2) | Out-Null
$d.Content.Find.Execute("81+16=97", $true, $false, $false, $false, $false, $true, 1, $false, "59-11=48", 2) | Out-Null
$d.Content.Find.Execute("13+80=93", $true, $false, $false, $false, $false, $true, 1, $false, "27+61=88", 2) | Out-Null
$d.Content.Find.Execute("65+17=82", $true, $false, $false, $false, $false, $true, 1, $false, "58+39=97", 2) | Out-Null
$d.Content.Find.Execute("38+12=50", $true, $false, $false, $false, $false, $true, 1, $false, "49-22=27", 2) | Out-Null
$d.Content.Find.Execute("51-19=32", $true, $false, $false, $false, $false, $true, 1, $false, "50+44=94", 2) | Out-Null
$d.Content.Find.Execute("5+80=85", $true, $false, $false, $false, $false, $true, 1, $false, "84+14=98", 2) | Out-Null
$d.Content.Find.Execute("57-27=30", $true, $false, $false, $false, $false, $true, 1, $false, "46+37=83", 2) | Out-Null
$d.Content.Find.Execute("89-11=78", $true, $false, $false, $false, $false, $true, 1, $false, "60-38=22", 2) | Out-Null
$d.Content.Find.Execute("40+51=91", $true, $false, $false, $false, $false, $true, 1, $false, "34+10=44", 2) | Out-Null
$d.Content.Find.Execute("15+71=86", $true, $false, $false, $false, $false, $true, 1, $false, "96-35=61", 2) | Out-Null
$d.Content.Find.Execute("38+3=41", $true, $false, $false, $false, $false, $true, 1, $false, "6+11=17", 2) | Out-Null
$d.Content.Find.Execute("61-54=7", $true, $false, $false, $false, $false, $true, 1, $false, "40+15=55", 2) | Out-Null
$d.Content.Find.Execute("51+10=61", $true, $false, $false, $false, $false, $true, 1, $false, "12+18=30", 2) | Out-Null
$d.Content.Find.Execute("97-53=44", $true, $false, $false, $false, $false, $true, 1, $false, "83-38=45", 2) | Out-Null
$d.Content.Find.Execute("53-2=51", $true, $false, $false, $false, $false, $true, 1, $false, "66-46=20", 2) | Out-Null
$d.Content.Find.Execute("16-15=1", $true, $false, $false, $false, $false, $true, 1, $false, "59+8=67", 2) | Out-Null
$d.Content.Find.Execute("67-60=7", $true, $false, $false, $false, $false, $true, 1, $false, "88-46=42", 2) | Out-Null
$d.Content.Find.Execute("70+29=99", $true, $false, $false, $false, $false, $true, 1, $false, "92-43=49", 2) | Out-Null
$d.Content.Find.Execute("4+7=11", $true, $false, $false, $false, $false, $true, 1, $false, "79-3=76", 2) | Out-Null
$d.Content.Find.Execute("33+56=89", $true, $false, $false, $false, $false, $true, 1, $false, "37+13=50", 2) | Out-Null
$d.Content.Find.Execute("90-88=2", $true, $false, $false, $false, $false, $true, 1, $false, "76-40=36", 2) | Out-Null
$d.Content.Find.Execute("12+14=26", $true, $false, $false, $false, $false, $true, 1, $false, "13+43=56", 2) | Out-Null
$d.Content.Find.Execute("68-26=42", $true, $false, $false, $false, $false, $true, 1, $false, "45+46=91", 2) | Out-Null
$d.Content.Find.Execute("97-84=13", $true, $false, $false, $false, $false, $true, 1, $false, "96-31=65", 2) | Out-Null
$d.Content.Find.Execute("20+65=85", $true, $false, $false, $false, $false, $true, 1, $false, "16+45=61", 2) | Out-Null
$d.Content.Find.Execute("98-81=17", $true, $false, $false, $false, $false, $true, 1, $false, "4+64=68", 2) | Out-Null
$d.Content.Find.Execute("3+47=50", $true, $false, $false, $false, $false, $true, 1, $false, "91-47=44", 2) | Out-Null
$d.Content.Find.Execute("45+33=78", $true, $false, $false, $false, $false, $true, 1, $false, "42+3=45", 2) | Out-Null
$d.Content.Find.Execute("85-71=14", $true, $false, $false, $false, $false, $true, 1, $false, "8+64=72", 2) | Out-Null
$d.Content.Find.Execute("39+5=44", $true, $false, $false, $false, $false, $true, 1, $false, "7+73=80", 2) | Out-Null
$d.Content.Find.Execute("49+1=50", $true, $false, $false, $false, $false, $true, 1, $false, "86-58=28", 2) | Out-Null
$d.Content.Find.Execute("38+1=39", $true, $false, $false, $false, $false, $true, 1, $false, "5+74=79", 2) | Out-Null
$d.Content.Find.Execute("59-57=2", $true, $false, $false, $false, $false, $true, 1, $false, "27+69=96", 2) | Out-Null
$d.Content.Find.Execute("32+35=67", $true, $false, $false, $false, $false, $true, 1, $false, "42+1=43", 2) | Out-Null
$d.Content.Find.Execute("93-80=13", $true, $false, $false, $false, $false, $true, 1, $false, "46+1=47", 2) | Out-Null
$d.Content.Find.Execute("73-37=36", $true, $false, $false, $false, $false, $true, 1, $false, "67-38=29", 2) | Out-Null
$d.Content.Find.Execute("62+15=77", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=81", 2) | Out-Null
$d.Content.Find.Execute("24+28=52", $true, $false, $false, $false, $false, $true, 1, $false, "89-54=35", 2) | Out-Null
$d.Content.Find.Execute("96-88=8", $true, $false, $false, $false, $false, $true, 1, $false, "31+39=70", 2) | Out-Null
$d.Content.Find.Execute("24+61=85", $true, $false, $false, $false, $false, $true, 1, $false, "21+55=76", 2) | Out-Null
$d.Content.Find.Execute("84-47=37", $true, $false, $false, $false, $false, $true, 1, $false, "32-25=7", 2) | Out-Null
$d.Content.Find.Execute("97-2=95", $true, $false, $false, $false, $false, $true, 1, $false, "37+42=79", 2) | Out-Null
$d.Content.Find.Execute("24+1=25", $true, $false, $false, $false, $false, $true, 1, $false, "0+35=35", 2) | Out-Null
$d.Content.Find.Execute("64-24=40", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=65", 2) | Out-Null
$d.Content.Find.Execute("46+47=93", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=77", 2) | Out-Null
$d.Content.Find.Execute("11+18=29", $true, $false, $false, $false, $false, $true, 1, $false, "93-19=74", 2) | Out-Null
$d.Content.Find.Execute("41-35=6", $true, $false, $false, $false, $false, $true, 1, $false, "39+45=84", 2) | Out-Null
$d.Content.Find.Execute("4+48=52", $true, $false, $false, $false, $false, $true, 1, $false, "21+59=80", 2) | Out-Null
$d.Content.Find.Execute("59-54=5", $true, $false, $false, $false, $false, $true, 1, $false, "41+11=52", 2) | Out-Null
$d.Content.Find.Execute("41-0=41", $true, $false, $false, $false, $false, $true, 1, $false, "10-4=6", 2) | Out-Null
$d.Content.Find.Execute("97-68=29", $true, $false, $false, $false, $false, $true, 1, $false, "63-32=31", 2) | Out-Null
$d.Content.Find.Execute("89-22=67", $true, $false, $false, $false, $false, $true, 1, $false, "84-63=21", 2) | Out-Null
$d.Content.Find.Execute("22+34=56", $true, $false, $false, $false, $false, $true, 1, $false, "44+20=64", 2) | Out-Null
$d.Content.Find.Execute("31+8=39", $true, $false, $false, $false, $false, $true, 1, $false, "90-81=9", 2) | Out-Null
$d.Content.Find.Execute("62-27=35", $true, $false, $false, $false, $false, $true, 1, $false, "67-53=14", 2) | Out-Null
$d.Content.Find.Execute("19+36=55", $true, $false, $false, $false, $false, $true, 1, $false, "69-7=62", 2) | Out-Null
$d.Content.Find.Execute("99-10=89", $true, $false, $false, $false, $false, $true, 1, $false, "81-74=7", 2) | Out-Null
$d.Content.Find.Execute("41+44=85", $true, $false, $false, $false, $false, $true, 1, $false, "11+65=76", 2) | Out-Null
$d.Content.Find.Execute("78-63=15", $true, $false, $false, $false, $false, $true, 1, $false, "58-14=44", 2) | Out-Null
$d.Content.Find.Execute("83-62=21", $true, $false, $false, $false, $false, $true, 1, $false, "21+50=71", 2) | Out-Null
$d.Content.Find.Execute("39-32=7", $true, $false, $false, $false, $false, $true, 1, $false, "41-36=5", 2) | Out-Null
$d.Content.Find.Execute("43+19=62", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=74", 2) | Out-Null
$d.Content.Find.Execute("3+21=24", $true, $false, $false, $false, $false, $true, 1, $false, "65-10=55", 2) | Out-Null
$d.Content.Find.Execute("46-22=24", $true, $false, $false, $false, $false, $true, 1, $false, "26-24=2", 2) | Out-Null
$d.Content.Find.Execute("29-16=13", $true, $false, $false, $false, $false, $true, 1, $false, "41+57=98", 2) | Out-Null
$d.Content.Find.Execute("34-21=13", $true, $false, $false, $false, $false, $true, 1, $false, "90-41=49", 2) | Out-Null
$d.Content.Find.Execute("86-62=24", $true, $false, $false, $false, $false, $true, 1, $false, "86-29=57", 2) | Out-Null
$d.Content.Find.Execute("78+10=88", $true, $false, $false, $false, $false, $true, 1, $false, "83-7=76", 2) | Out-Null
$d.Content.Find.Execute("63-31=32", $true, $false, $false, $false, $false, $true, 1, $false, "20+38=58", 2) | Out-Null
$d.Content.Find.Execute("41+38=79", $true, $false, $false, $false, $false, $true, 1, $false, "47-29=18", 2) | Out-Null
$d.Content.Find.Execute("94+1=95", $true, $false, $false, $false, $false, $true, 1, $false, "15+50=65", 2) | Out-Null
$d.Content.Find.Execute("79+19=98", $true, $false, $false, $false, $false, $true, 1, $false, "85-44=41", 2) | Out-Null
$d.Content.Find.Execute("14+15=29", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=60", 2) | Out-Null
$d.Content.Find.Execute("45-40=5", $true, $false, $false, $false, $false, $true, 1, $false, "6+50=56", 2) | Out-Null
$d.Content.Find.Execute("86-14=72", $true, $false, $false, $false, $false, $true, 1, $false, "72-27=45", 2) | Out-Null
$d.Content.Find.Execute("80-73=7", $true, $false, $false, $false, $false, $true, 1, $false, "22+39=61", 2) | Out-Null
$d.Content.Find.Execute("44+35=79", $true, $false, $false, $false, $false, $true, 1, $false, "62-32=30", 2) | Out-Null
$d.Content.Find.Execute("7+36=43", $true, $false, $false, $false, $false, $true, 1, $false, "99-88=11", 2) | Out-Null
$d.Content.Find.Execute("91-61=30", $true, $false, $false, $false, $false, $true, 1, $false, "74-46=28", 2) | Out-Null
$d.Content.Find.Execute("96-32=64", $true, $false, $false, $false, $false, $true, 1, $false, "65-55=10", 2) | Out-Null
$d.Content.Find.Execute("59+16=75", $true, $false, $false, $false, $false, $true, 1, $false, "24+35=59", 2) | Out-Null
$d.Content.Find.Execute("99-66=33", $true, $false, $false, $false, $false, $true, 1, $false, "27+16=43", 2) | Out-Null
$d.Content.Find.Execute("63-51=12", $true, $false, $false, $false, $false, $true, 1, $false, "40+52=92", 2) | Out-Null
$d.Content.Find.Execute("6+75=81", $true, $false, $false, $false, $false, $true, 1, $false, "49+8=57", 2) | Out-Null
$d.Content.Find.Execute("19+61=80", $true, $false, $false, $false, $false, $true, 1, $false, "61+26=87", 2) | Out-Null
$d.Content.Find.Execute("92-17=75", $true, $false, $false, $false, $false, $true, 1, $false, "43-3=40", 2) | Out-Null
$d.Content.Find.Execute("24-20=4", $true, $false, $false, $false, $false, $true, 1, $false, "52-16=36", 2) | Out-Null
$d.Content.Find.Execute("53-20=33", $true, $false, $false, $false, $false, $true, 1, $false, "99-14=85", 2) | Out-Null
$d.Content.Find.Execute("39+40=79", $true, $false, $false, $false, $false, $true, 1, $false, "53-14=39", 2) | Out-Null
$d.Content.Find.Execute("40+23=63", $true, $false, $false, $false, $false, $true, 1, $false, "24+59=83", 2) | Out-Null
$d.Content.Find.Execute("0+52=52", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=18", 2) | Out-Null
